$d = $word.ActiveDocument

# En dash character used in the heading ("Plano de Atendimento – PDA")
$enDash = [char]0x2013
$target = "Plano de Atendimento " + $enDash + " PDA"

# Locate the exact run of text we need to split (the top-right title,
# NOT the other headings like "Objetivo do Plano de Atendimento - PDA").
$searchRange = $d.Content
$found = $searchRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target heading text"
}

# Re-create a plain Range from the found boundaries: re-using the Find
# range object directly with InsertXML leaves the first couple of
# characters of the old text behind, so work off a fresh Range instead.
$findRange = $d.Range($searchRange.Start, $searchRange.End)

# Build the OOXML fragment that replaces the single run with three runs
# plus spell-check proof-error markers around "PDAa", matching the
# target diff: "Plano de Atendimento - " + ("PDA" / "a" wrapped in
# proofErr spellStart/spellEnd).
$openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body><w:p>' +
'<w:r w:rsidRPr="00FB6CCE"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Plano de Atendimento ' + $enDash + ' </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>PDA</w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>a</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'</w:p></w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$findRange.InsertXML($openXml)
